$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# This sheet is a "from -> to" field-rename map, grouped by section
# (documents / studies / subjects / series / conc_time_values), used to
# convert the old cvtdb column names into the new template's column names.
# The edit adds several newly-introduced fields in the various sections and
# renames two "studies" fields to their "_original" suffixed variants.
#
# We insert the new rows from the bottom of the sheet upward so that the
# row numbers used for the still-untouched, earlier insert points remain
# valid (no need to re-derive shifted row numbers as we go).
# ---------------------------------------------------------------------------

# 1) "series" section: new "conc_cumulative" row, right after "y_max" (old row 57)
$ws.Rows("58:58").Insert()
$ws.Range("A58").Value = "series"
$ws.Range("B58").Value = "conc_cumulative"
$ws.Range("C58").Value = "conc_cumulative"

# 2) "studies" section: new "dose_volume_units" row, right after "dose_volume" (old row 22)
$ws.Rows("23:23").Insert()
$ws.Range("A23").Value = "studies"
$ws.Range("B23").Value = "dose_volume_units"
$ws.Range("C23").Value = "dose_volume_units"

# 3) "studies" section: two existing rows get their "from" column renamed to the
#    "_original" suffixed source field (the "to" / template column stays the same)
$ws.Range("B12").Value = "test_substance_name_secondary_original"
$ws.Range("B13").Value = "test_substance_casrn_original"

# 4) "studies" section: new "fk_reference_document_id" row, right before "id" (old row 10)
$ws.Rows("10:10").Insert()
$ws.Range("A10").Value = "studies"
$ws.Range("B10").Value = "fk_reference_document_id"
$ws.Range("C10").Value = "fk_reference_document_id"

# 5) "documents" section: new "clowder_file_id" row, appended at the end of the
#    section (right before "studies" begins at old row 10)
$ws.Rows("10:10").Insert()
$ws.Range("A10").Value = "documents"
$ws.Range("B10").Value = "clowder_file_id"
$ws.Range("C10").Value = "clowder_file_id"

# 6) "documents" section: three new rows inserted at the top of the section
#    (right after the header row)
$ws.Rows("2:2").Insert()
$ws.Range("A2").Value = "documents"
$ws.Range("B2").Value = "id"
$ws.Range("C2").Value = "id"

$ws.Rows("3:3").Insert()
$ws.Range("A3").Value = "documents"
$ws.Range("B3").Value = "document_type"
$ws.Range("C3").Value = "document_type"

$ws.Rows("4:4").Insert()
$ws.Range("A4").Value = "documents"
$ws.Range("B4").Value = "extracted"
$ws.Range("C4").Value = "extracted"

# ---------------------------------------------------------------------------
# Refresh the autofilter / defined name / dimension so they cover the grown
# range A1:C85 (sheet grew from 78 to 85 data+header rows).
# ---------------------------------------------------------------------------
$ws.AutoFilterMode = $false
$ws.Range("A1:C85").AutoFilter()

foreach ($n in $wb.Names) {
    if ($n.Name -eq "Sheet1!_FilterDatabase") {
        $n.RefersTo = "=Sheet1!`$A`$1:`$C`$85"
    }
}

# Restore the selection to where the author last left the cursor.
$ws.Range("C13").Select()
